# Generate Report for Handoff
# Updates the localization-status report: Status flips from
# "In Translation" to "Ready for handoff" and the associated timestamps
# move forward to the handoff generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2016-08-31 05:00:14"
$dede.Range("H2").Value = "2016-08-31 05:00:14"
$zhcn.Range("H2").Value = "2016-08-31 04:59:57"

# --- Column widths: the longer "Ready for handoff" label needs a wider
#     Status column on every sheet that shows it. ---
$overview.Columns.Item(5).ColumnWidth = 16.3826517722958
$overview.Columns.Item(6).ColumnWidth = 16.3826517722958
$zhcn.Columns.Item(3).ColumnWidth = 16.3826517722958
$dede.Columns.Item(3).ColumnWidth = 16.3826517722958
